$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new triplicate titration result row (2021-04-28) below row 30.
$ws.Cells.Item(31, 1).Value = 20210428
$ws.Cells.Item(31, 2).Value = 2223.7829999999999
$ws.Cells.Item(31, 3).Value = 2225.4699999999998
$ws.Cells.Item(31, 4).Formula = "=100*(B31-C31)/C31"
$ws.Cells.Item(31, 5).Value = 181
$ws.Cells.Item(31, 6).Value = $ws.Cells.Item(30, 6).Value2

# Keep the selection/view in sync with where Excel would land after data entry.
$ws.Range("F32").Select()
